# Append a date run (" 1/6/23") right after the "Meeting notes" title text,
# keeping it as its own run with the same character formatting
# (accent6 green color, 32 half-point size, en-GB language) as the title.

$d = $word.ActiveDocument

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

# Sanity check: make sure we are editing the expected paragraph.
if ($titleRange.Text -notmatch "^Meeting notes") {
    throw "Unexpected first paragraph text: '$($titleRange.Text)'"
}

# Collapse to just before the paragraph mark so the new run lands inside
# the same paragraph, after the existing "Meeting notes" run.
$insertionPoint = $d.Range($titleRange.Start, $titleRange.End - 1)

$openXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"
                   xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="2C3ECCC9" w14:textId="652A97FD" w:rsidR="00463D63" w:rsidRPr="00C16F74" w:rsidRDefault="00463D63" w:rsidP="00C16F74">
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:color w:val="70AD47" w:themeColor="accent6"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="00C16F74">
              <w:rPr>
                <w:color w:val="70AD47" w:themeColor="accent6"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t>Meeting notes</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:color w:val="70AD47" w:themeColor="accent6"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t xml:space="preserve"> 1/6/23</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$titleRange.InsertXML($openXml)
